# Ferguson_URL.xlsx update: rebuild the "URL" sheet with the new
# Type / Short Name columns, refreshed Homologo Mansfield data and a
# single remaining hyperlink on the last row, per "Update info capture
# from master_url file".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlPasteValues  = -4163
$xlLeft         = -4131

# ---------------------------------------------------------------------
# Phase 1: harvest the existing cell formats (from the ORIGINAL A1:J5
# layout) into an out-of-the-way palette row before anything is
# overwritten. Each palette cell ends up carrying exactly one of the
# distinct styles used by the rebuilt table.
# ---------------------------------------------------------------------
$paletteRow = 500

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

Copy-Format "D1" "A$paletteRow"    # style 1
Copy-Format "E1" "B$paletteRow"    # style 2
Copy-Format "A2" "C$paletteRow"    # style 3
Copy-Format "D2" "D$paletteRow"    # style 4
Copy-Format "E3" "E$paletteRow"    # style 5
Copy-Format "F2" "F$paletteRow"    # style 6
Copy-Format "J2" "G$paletteRow"    # style 7
Copy-Format "E2" "H$paletteRow"    # style 8
Copy-Format "H2" "I$paletteRow"    # style 9
Copy-Format "H3" "J$paletteRow"    # style 10

# style 11 is new: same border as style 4 (D2) but left-aligned instead
# of centered.
Copy-Format "D2" "K$paletteRow"
$ws.Range("K$paletteRow").HorizontalAlignment = $xlLeft

$excel.CutCopyMode = $false

$pA = "A$paletteRow"; $pB = "B$paletteRow"; $pC = "C$paletteRow"; $pD = "D$paletteRow"
$pE = "E$paletteRow"; $pF = "F$paletteRow"; $pG = "G$paletteRow"; $pH = "H$paletteRow"
$pI = "I$paletteRow"; $pJ = "J$paletteRow"; $pK = "K$paletteRow"

# ---------------------------------------------------------------------
# Phase 2: wipe the old hyperlinks and the old table body.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Range("A1:J5").Clear()

# ---------------------------------------------------------------------
# Phase 3: helper to write a cell's value + formatting in one go.
# $forceText ($true/$false) forces numeric-looking strings (SKU-style
# digit codes) to stay text instead of being auto-coerced to numbers -
# it round-trips the value through a scratch cell using the leading
# apostrophe trick, then pastes VALUES ONLY so the target cell's style
# is untouched by the quote-prefix flag that trick leaves behind.
# ---------------------------------------------------------------------
function Set-Cell($addr, $value, $paletteAddr, $forceText) {
    if ($forceText) {
        $ws.Range("ZZ1").Value = "'" + $value
        $ws.Range("ZZ1").Copy()
        $ws.Range($addr).PasteSpecial($xlPasteValues)
    } else {
        $ws.Range($addr).Value = $value
    }
    $ws.Range($paletteAddr).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# Row 1 - headers
Set-Cell "A1" "Fabricante"         $pB
Set-Cell "B1" "Homologo Mansfield" $pB
Set-Cell "C1" "Sku"                $pB
Set-Cell "D1" "Type"               $pA
Set-Cell "E1" "Linea"              $pA
Set-Cell "F1" "Rough in"           $pB
Set-Cell "G1" "Bowl Height"        $pB
Set-Cell "H1" "Asiento"            $pB
Set-Cell "I1" "Capacidad (Gpl)"    $pB
Set-Cell "J1" "Descripcion"        $pB
Set-Cell "K1" "Short Name"         $pB
Set-Cell "L1" "Link"               $pB

# Row 2 - Gerber Maxwell ADA EL Bowl
Set-Cell "A2" "Gerber"             $pC
Set-Cell "B2" "5916CTK"            $pC
Set-Cell "C2" "G0021975"           $pC
Set-Cell "D2" "Bowl"               $pK
Set-Cell "E2" "Maxwell"            $pD
Set-Cell "F2" "4 - 1/4"            $pH
Set-Cell "G2" "ADA"                $pF
Set-Cell "H2" "Elongated"          $pF
Set-Cell "I2" 1.28                 $pI
Set-Cell "J2" "Maxwell® 1.28 gpf Elongated Floor Mount Two Piece Toilet Bowl in White" $pF
Set-Cell "K2" "Gerber Maxwell ADA EL Bowl" $pF
Set-Cell "L2" "https://www.ferguson.com/product/gerber-plumbing-maxwell-128-gpf-elongated-floor-mount-two-piece-toilet-bowl-in-white-gg0021975/_/R-4493463" $pG

# Row 3 - Gerber Maxwell 1,28 gpf Tank
Set-Cell "A3" "Gerber"             $pC
Set-Cell "B3" "317310000"          $pC
Set-Cell "C3" "G0028990"           $pC
Set-Cell "D3" "Tank"               $pK
Set-Cell "E3" "Maxwell"            $pD
Set-Cell "F3" "NA"                 $pE
Set-Cell "G3" "NA"                 $pC
Set-Cell "H3" "NA"                 $pC
Set-Cell "I3" 1.28                 $pJ
Set-Cell "J3" "Maxwell® 1.28 gpf Toilet Tank with Left-Hand Trip Lever in White" $pF
Set-Cell "K3" "Gerber Maxwell 1,28 gpf Tank" $pF
Set-Cell "L3" "https://www.ferguson.com/product/gerber-plumbing-maxwell-128-gpf-toilet-tank-with-left-hand-trip-lever-in-white-gg0028990/_/R-4239960" $pG

# Row 4 - Gerber Maxwell ADA EL Bowl (different SKU)
Set-Cell "A4" "Gerber"             $pC
Set-Cell "B4" "137210040"          $pC
Set-Cell "C4" "GMX21928"           $pC
Set-Cell "D4" "Bowl"               $pK
Set-Cell "E4" "Maxwell"            $pD
Set-Cell "F4" "14"                 $pE -ForceText
Set-Cell "G4" "ADA"                $pC
Set-Cell "H4" "Elongated"          $pC
Set-Cell "I4" 1.28                 $pJ
Set-Cell "J4" "Maxwell® 1.28 gpf Elongated Toilet Bowl in White" $pF
Set-Cell "K4" "Gerber Maxwell ADA EL Bowl" $pF
Set-Cell "L4" "https://www.ferguson.com/product/gerber-plumbing-maxwell-128-gpf-elongated-toilet-bowl-in-white-ggmx21928/_/R-7562151?trackSignal=true" $pG

# Row 5 - Gerber Maxwell 1,6 gpf Tank
Set-Cell "A5" "Gerber"             $pC
Set-Cell "B5" "160010007"          $pC
Set-Cell "C5" "GMX28990"           $pC
Set-Cell "D5" "Tank"               $pK
Set-Cell "E5" "Maxwell"            $pD
Set-Cell "F5" "NA"                 $pE
Set-Cell "G5" "NA"                 $pC
Set-Cell "H5" "NA"                 $pC
Set-Cell "I5" 1.6                  $pJ
Set-Cell "J5" "Maxwell® 1.6 gpf Two Piece Toilet Tank in White" $pF
Set-Cell "K5" "Gerber Maxwell 1,6 gpf Tank" $pF
Set-Cell "L5" "https://www.ferguson.com/product/gerber-plumbing-maxwell-16-gpf-two-piece-toilet-tank-in-white-ggmx28990/_/R-7005560" $pG

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Phase 4: the only surviving hyperlink (Gerber Maxwell 1,6 gpf Tank).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("L5"), "https://www.ferguson.com/product/gerber-plumbing-maxwell-16-gpf-two-piece-toilet-tank-in-white-ggmx28990/_/R-7005560") | Out-Null

# ---------------------------------------------------------------------
# Phase 5: clean up the scratch cells used while building the table.
# ---------------------------------------------------------------------
$ws.Range("ZZ1").Clear()
$ws.Range("A$($paletteRow):K$($paletteRow)").Clear()

# ---------------------------------------------------------------------
# Phase 6: column widths for the now-12-column layout.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth  = 17.42578125
$ws.Columns.Item(2).ColumnWidth  = 26.28515625
$ws.Columns.Item(3).ColumnWidth  = 13.7109375
$ws.Columns.Item(4).ColumnWidth  = 13.7109375
$ws.Columns.Item(5).ColumnWidth  = 14.42578125
$ws.Columns.Item(6).ColumnWidth  = 16.28515625
$ws.Columns.Item(7).ColumnWidth  = 16.42578125
$ws.Columns.Item(8).ColumnWidth  = 20.7109375
$ws.Columns.Item(9).ColumnWidth  = 20.7109375
$ws.Columns.Item(10).ColumnWidth = 70.7109375
$ws.Columns.Item(11).ColumnWidth = 29.28515625
$ws.Columns.Item(12).ColumnWidth = 139.7109375

# ---------------------------------------------------------------------
# Phase 7: selection / scroll position, matching the saved view.
# ---------------------------------------------------------------------
$ws.Range("L15").Select()
